$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Currently:
#   Row17 = 1043644876 / ROBERT ENRIQUE PAZOS ALANDETE / 2204 / 32000 / 1000000
#   Row18 = 1044930744 / ERNEIS JOHED PAJARO CASTRO     / 2206 / 40000 / 1200000
#   Row19 = 1047435771 / GUSTAVO ADOLFO MARTINEZ PINTO  / 2203 / 21333 / 1000000
#
# New record (GUSTAVO) is inserted right after row 16, pushing ROBERT and
# ERNEIS down by one row each. Net effect: a rotation of the three data rows.

$ws.Range("C17").Value = "1047435771"
$ws.Range("D17").Value = "GUSTAVO ADOLFO MARTINEZ PINTO"
$ws.Range("E17").Value = "2203"
$ws.Range("F17").Value = 21333
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "1043644876"
$ws.Range("D18").Value = "ROBERT ENRIQUE PAZOS ALANDETE"
$ws.Range("E18").Value = "2204"
$ws.Range("F18").Value = 32000
$ws.Range("G18").Value = 1000000

$ws.Range("C19").Value = "1044930744"
$ws.Range("D19").Value = "ERNEIS JOHED PAJARO CASTRO"
$ws.Range("E19").Value = "2206"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1200000
